# Apply "Add data for 2022-12-06" update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (and matching "through" label) from 11-27 to 11-28
$ws.Name = "Through 2022-11-28"
$ws.Range("I1").Value = "2022 (through 11-28)"

# Update November (row 12) and Total (row 14) counts for the 2022 column (I)
$ws.Range("I12").Value = 106
$ws.Range("I14").Value = 1504
